# Natmi following Dr Hou advice
# Re-derive the Bsg -> Slc16a7 ligand/receptor table using 3 senders x 3 targets
# (ECs, FAPs, sCs) with the updated "Ligand-expressing cells" count of 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @("ECs",  "Bsg", "Slc16a7", "ECs",  3, 1, 49.400308,          148.200924,         0.3028101582105581, 0.3028101582105581, 2, 0.6666666666666666, 1.387792333333333, 4.163377, 0.4491959817916776, 0.4491959817916776, 68.55736870670532, 617.0163183603479, 0.1360211063138849,  0.1360211063138848),
    @("ECs",  "Bsg", "Slc16a7", "FAPs", 3, 1, 49.400308,          148.200924,         0.3028101582105581, 0.3028101582105581, 2, 0.6666666666666666, 0.4615953333333334, 1.384786, 0.1494076339570906, 0.1494076339570906, 22.80295163802933, 205.226564742264,  0.04524214927641176, 0.04524214927641176),
    @("ECs",  "Bsg", "Slc16a7", "sCs",  3, 1, 49.400308,          148.200924,         0.3028101582105581, 0.3028101582105581, 3, 1,                  1.240115333333333, 3.720346, 0.4013963842512318, 0.4013963842512318, 61.26207942218932, 551.3587147997039, 0.1215469026202615,  0.1215469026202615),
    @("FAPs", "Bsg", "Slc16a7", "ECs",  3, 1, 69.564149,          208.692447,         0.4264088994034781, 0.4264088994034782, 2, 0.6666666666666666, 1.387792333333333, 4.163377, 0.4491959817916776, 0.4491959817916776, 96.54059265705766, 868.865333913519,  0.1915411642122541,  0.1915411642122541),
    @("FAPs", "Bsg", "Slc16a7", "FAPs", 3, 1, 69.564149,          208.692447,         0.4264088994034781, 0.4264088994034782, 2, 0.6666666666666666, 0.4615953333333334, 1.384786, 0.1494076339570906, 0.1494076339570906, 32.11048654570467, 288.994378911342,  0.06370874475812074, 0.06370874475812074),
    @("FAPs", "Bsg", "Slc16a7", "sCs",  3, 1, 69.564149,          208.692447,         0.4264088994034781, 0.4264088994034782, 3, 1,                  1.240115333333333, 3.720346, 0.4013963842512318, 0.4013963842512318, 86.26756782518466, 776.408110426662,  0.1711589904331034,  0.1711589904331034),
    @("sCs",  "Bsg", "Slc16a7", "ECs",  3, 1, 44.17507666666666,  132.52523,          0.2707809423859638, 0.2707809423859638, 2, 0.6666666666666666, 1.387792333333333, 4.163377, 0.4491959817916776, 0.4491959817916776, 61.30583272241221, 551.75249450171,   0.1216337112655387,  0.1216337112655387),
    @("sCs",  "Bsg", "Slc16a7", "FAPs", 3, 1, 44.17507666666666,  132.52523,          0.2707809423859638, 0.2707809423859638, 2, 0.6666666666666666, 0.4615953333333334, 1.384786, 0.1494076339570906, 0.1494076339570906, 20.39100923897556, 183.51908315078,   0.04045673992255813, 0.04045673992255812),
    @("sCs",  "Bsg", "Slc16a7", "sCs",  3, 1, 44.17507666666666,  132.52523,          0.2707809423859638, 0.2707809423859638, 3, 1,                  1.240115333333333, 3.720346, 0.4013963842512318, 0.4013963842512318, 54.78218992550888, 493.0397093295799, 0.108690491197867,   0.108690491197867)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $excelRow = $r + 2
    $rowVals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $excelRow).Value = $rowVals[$i]
    }
}
